# Se procesan de nuevo los datos con las nuevas dimensiones curadas
#
# Column D (corine-land-cover-2000-nivel-3-descripcion) moves from being a
# curated "dimension" to a "measure", while column F (municipio-nombre)
# moves from being a "measure" to a curated "dimension" (refArea), gaining
# a URI-Municipio mapping column. The now-unused mapping file reference for
# column D (row 5) is cleared since it is no longer a dimension.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: iaest-dimension/iaest-measure annotation row
$ws.Range("D2").Value = "iaest-measure:corine-land-cover-2000-nivel-3-descripcion"
$ws.Range("F2").Value = "sdmx-dimension:refArea"

# Row 3: dim / medida row
$ws.Range("D3").Value = "medida"
$ws.Range("F3").Value = "dim"

# Row 4: type / URI-mapping-column row
$ws.Range("D4").Value = "xsd:int"
$ws.Range("F4").Value = "URI-Municipio"

# Row 5: mapping workbook filename row - column D no longer needs a mapping
# file, since it is now a plain measure instead of a curated dimension.
$ws.Range("D5").Value = ""
